$d = $word.ActiveDocument

# 1. "ithimbeni" -> "ithimba" (in the chatbot/email sentence)
$d.Content.Find.Execute(
    "ithimbeni locwaningo ku-", $true, $false, $false, $false, $false,
    $true, 1, $false, "ithimba locwaningo ku-", 2)

# 2. "uzizwa" -> "uzizwe" (in the safe-space sentence)
$d.Content.Find.Execute(
    "uzizwa ukhululekile uma", $true, $false, $false, $false, $false,
    $true, 1, $false, "uzizwe ukhululekile uma", 2)

# 3. "kanti iMenenja yocwaningo nguZamakhanya" -> "kanye neMenenja yocwaningo uZamakhanya"
$d.Content.Find.Execute(
    "kanti iMenenja yocwaningo nguZamakhanya", $true, $false, $false, $false, $false,
    $true, 1, $false, "kanye neMenenja yocwaningo uZamakhanya", 2)

# 4. "noma okukukhathazayo" -> "noma kukhona okukukhathazayo" (rights sentence)
$d.Content.Find.Execute(
    "Uma unemibuzo noma okukukhathazayo mayelana", $true, $false, $false, $false, $false,
    $true, 1, $false, "Uma unemibuzo noma kukhona okukukhathazayo mayelana", 2)

# 5. "idokhumenti" -> "incwadi" (WhatsApp consent sentence)
$d.Content.Find.Execute(
    "waqonda idokhumenti engenhla", $true, $false, $false, $false, $false,
    $true, 1, $false, "waqonda incwadi engenhla", 2)
